$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new data row above the current totals row (row 29) and seed it
#    by duplicating the formatting / merged-cell layout of the last existing
#    data row (row 28). This keeps styles, borders and merges identical to
#    the other item rows without touching far-away columns.
# ---------------------------------------------------------------------------
$ws.Rows.Item(29).Insert()
$ws.Range("A28:Q28").Copy($ws.Range("A29:Q29"))
$ws.Rows.Item(29).RowHeight = 25.5

# Helper: assign a text value to a cell while forcing it to be stored as a
# genuine text/shared-string, even when the cell's display number format is
# numeric (otherwise Excel would silently convert a numeric-looking string
# such as "8.0000" into a plain number and drop the t="s" flag).
function Set-TextValue($range, [string]$text) {
    $originalFormat = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = $originalFormat
}

# ---------------------------------------------------------------------------
# 2. Populate the new item row (#23 "مجموعه برد").
# ---------------------------------------------------------------------------
$ws.Range("A29").Value = 23
Set-TextValue $ws.Range("C29") "مجموعه برد"
Set-TextValue $ws.Range("H29") "0:0"
Set-TextValue $ws.Range("N29") "8.00"
Set-TextValue $ws.Range("P29") "8.0000"
# L29 ("0") and Q29 ("1:0") already match the required values after the copy
# from row 28, so they do not need to be touched.

# ---------------------------------------------------------------------------
# 3. Update the totals row (now row 30): new grand total and row height.
# ---------------------------------------------------------------------------
$ws.Range("P30").Value = 923.83
$ws.Rows.Item(30).RowHeight = 24.75

# ---------------------------------------------------------------------------
# 4. Update the footer timestamp (now row 31).
# ---------------------------------------------------------------------------
$ws.Range("A31").Value = "Sunday, 22 June, 2025 11:15 AM"
